# Update dashboards - 2026-02-04
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: ADP Total NonFarm Private (date rolled forward a month; also gets
#     the "current" yellow-highlight format used by other up-to-date rows) ---
$ws.Range("N13").Copy()
$ws.Range("N5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("N5").Value = 46023

$ws.Range("Q5").Value = 22000
$ws.Range("R5").Value = 37000
$ws.Range("S5").Value = 74000
$ws.Range("T5").Value = 20000
$ws.Range("U5").Value = 88000

# --- Row 29: 5yr, 5yr Forward inflation ---
$ws.Range("N29").Value = 46056
$ws.Range("Q29").Value = 2.19
$ws.Range("R29").Value = 2.18
$ws.Range("S29").Value = 2.19
$ws.Range("T29").Value = 2.18
$ws.Range("U29").Value = 2.22

# --- Row 30: 10yr TIPS ---
$ws.Range("N30").Value = 46056
$ws.Range("Q30").Value = 2.36
$ws.Range("R30").Value = 2.35
$ws.Range("S30").Value = 2.36
$ws.Range("T30").Value = 2.35
$ws.Range("U30").Value = 2.36

# --- Row 47: FFR (date only) ---
$ws.Range("N47").Value = 46055

# --- Row 48: 2y UST ---
$ws.Range("N48").Value = 46055
$ws.Range("Q48").Value = 3.57
$ws.Range("R48").Value = 3.52
$ws.Range("S48").Value = 3.53
$ws.Range("T48").Value = 3.56
$ws.Range("U48").Value = 3.53

# --- Row 49: 5y UST ---
$ws.Range("N49").Value = 46055
$ws.Range("Q49").Value = 3.83
$ws.Range("R49").Value = 3.79
$ws.Range("S49").Value = 3.8
$ws.Range("T49").Value = 3.83
$ws.Range("U49").Value = 3.81

# --- Row 50: 10y UST ---
$ws.Range("N50").Value = 46055
$ws.Range("Q50").Value = 4.29
$ws.Range("R50").Value = 4.26
$ws.Range("S50").Value = 4.24
$ws.Range("T50").Value = 4.26
$ws.Range("U50").Value = 4.24

# --- Row 52: BAA ---
$ws.Range("N52").Value = 46055
$ws.Range("Q52").Value = 5.9
$ws.Range("R52").Value = 5.88
$ws.Range("S52").Value = 5.87
$ws.Range("T52").Value = 5.88
$ws.Range("U52").Value = 5.85
